# Minor changes to the "side menu collapsed" results sheet:
#  - rename the sheet from "Results" to "verifyingSideMenuCollapsedLangA"
#  - scroll the sheet view down so row 53 is the top-left visible row
#  - move the selection from L8 to G57

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "verifyingSideMenuCollapsedLangA"

# Make sure the sheet is the active one, then scroll so row 53 becomes the
# top-left visible cell of the window (best-effort - scrolls the window
# before the final selection is made below).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 53
$excel.ActiveWindow.ScrollColumn = 1

# Move the active selection to G57.
$ws.Range("G57").Select()
